# Generate Report for Handoff
# - Flip "Status" / latest-handoff cells from the stale "Handed back: in sync
#   with en-US" text to "Ready for handoff".
# - Refresh the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
#   timestamps to the new handoff run time.
# - Columns that held the (now much shorter) status text narrow to fit.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# ---------------------------------------------------------------------------
# Overview sheet: zh-cn / de-de status columns (E, F) + HO Xliff generate
# date (G).
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = "2016-09-03 05:05:03"

# ---------------------------------------------------------------------------
# zh-cn sheet: Status (C) + Latest Handoff Datetime (H).
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = "2016-09-03 05:04:57"

# ---------------------------------------------------------------------------
# de-de sheet: Status (C) + Latest Handoff Datetime (H).
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = "2016-09-03 05:05:03"

# ---------------------------------------------------------------------------
# Narrow the columns that used to hold the long status text so they fit the
# shorter "Ready for handoff" string.
# ---------------------------------------------------------------------------
$narrowWidth = 16.333333333333

$wsOverview.Columns.Item(5).ColumnWidth = $narrowWidth
$wsOverview.Columns.Item(6).ColumnWidth = $narrowWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $narrowWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $narrowWidth
